$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("U41_01", 41, "Nhiều mây", "Cloudy", "We are having a cloudy day today.", "cloudy day / ngày nhiều mây", "Adj"),
    @("U41_02", 41, "Làm bừng sáng", "Brighten up", "Your smile brightens up my day", "brighten up somebody's day / làm bừng sáng một ngày của ai đó", "V"),
    @("U41_03", 41, "Dự báo", "Predict", "No one can predict the future", "predict the future / dự đoán trước tương lai", "V"),
    @("U41_04", 41, "Chính xác", "Accurate", "Was what you said accurate information?", "accurate information", "Adj"),
    @("U41_05", 41, "Ôn hòa", "Mild", "Dalat has such as mild climate", "a mild climate / khí hậu ôn hòa", "Adj"),
    @("U41_06", 41, "Bão", "Storm", "A tropical storm swept through the village.", "A tropical storm ", "N"),
    @("U41_07", 41, "Mưa nặng hạt", "Pour down", "Suddenly, the rain pours down", "rain pours down / mưa rơi tầm tã", "V"),
    @("U41_08", 41, "Chớp", "Lightning", "A bolt of lightning struck down the house", "a bolt of lightning / một tia sét", "N"),
    @("U41_09", 41, "Sấm", "Thunder", "They suddenly heard a clap of thunder", "a clap of thunder / một tiếng sấm", "N"),
    @("U41_10", 41, "Sườn, dốc", "Slope", "We slid down the hill on a slope / Chúng tôi trượt xuống sườn đồi trên một con dốc", "on a slope / trên một dốc", "N"),
    @("U41_11", 41, "Nhiều nắng", "Sunny", "Let's not waste this sunny day / đừng lãng phí một ngày đầy nắng như thế này.", "a sunny day", "Adj"),
    @("U41_12", 41, "Mùa xuân", "Spring", "Families do spring cleaning before Te holiday.", "spring cleaning / tổng vệ sinh đầu xuân", "N"),
    @("U41_13", 41, "Bất thường", "unusual", "She behaved in a highly unusual manner", "highly unusual / rất bất thường", "Adj"),
    @("U41_14", 41, "Nơi trú ẩn", "Shelter", "You can volunteer at the homeless shelter", "a homeless shelter / chỗ ở cho người vô gia cư", "N"),
    @("U41_15", 41, "Mùa", "Season", "It's growing season for berries (quả mọng)", "growing season / mùa gieo trồng", "N"),
    @("U41_16", 41, "Sa mạc", "Desert", "It is hot like the sahara desert", "the Sahara Desert / sa mạc Sahara", "N"),
    @("U41_17", 41, "Mùa thu", "Fall", "Fall weather feels cooler than summer weather", "fall weather", "N"),
    @("U41_18", 41, "Cảnh", "Sight", "A traffic jam is a common sight in big cities", "common sight / cảnh thường thấy", "N"),
    @("U41_19", 41, "Địa lý", "Geography", "A degree in geography helps you understand the Earth", "a degree in geography / một tấm bằng ngành địa lý", "N"),
    @("U41_20", 41, "Nghiêm trọng, dữ dội", "Severe", "I took a pull for my severe headache", "a severe headache / một cơn đau đầu dữ dội", "Adj"),
    @("U41_21", 41, "Động đất", "Earthquake", "Volcanoes can cause a powerful earthquake.", "powerful earthquake / trận động đất mạnh", "N"),
    @("U41_22", 41, "Mặt đất", "Ground", "The ball is on the ground", "on the ground", "N"),
    @("U41_23", 41, "Trái đất", "Earth", "The earth is round", "the earth", "N"),
    @("U41_24", 41, "Tai họa", "Disaster", "the flood was a natural disaster", "a natural disaster / thiên tai", "N"),
    @("U41_25", 41, "Sóng", "Wave", "Surfers love to ride the wave", "ride the wave / cưỡi sóng", "N"),
    @("U41_26", 41, "Bờ biển", "Coast", "I have a house on the coast", "on the coast / trên bờ biển", "N"),
    @("U41_27", 41, "Đóng băng", "Freeze", "You can freeze fruits to make smoothies later", "freeze something", "V"),
    @("U41_28", 41, "Hồ", "Lake", "The campfire is by the lake", "by the lake / Bên hồ", "N"),
    @("U41_29", 41, "Phá hủy", "Destroy", "My niece (cháu gái) destroyed my toys", "destroy something", "V"),
    @("U41_30", 41, "Vụ mùa", "Crop", "We are short of food crops", "food crop / cây lương thực", "N")
)

$startRow = 1202
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $r = $data[$i]
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = [int]$r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
}

$ws.Range("C1186").Select() | Out-Null